$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header strings: "_old" suffix -> "_FV2410", "_new" suffix -> "_FV2504"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $orig = $cell.Value2
    if ($orig -ne $null) {
        $renamed = $orig.Replace("_old", "_FV2410").Replace("_new", "_FV2504")
        if ($renamed -ne $orig) {
            $cell.Value = $renamed
        }
    }
}

# 2. Freeze the header row (row 1) and keep selection on bottom pane
$ws.Activate()
[void]$ws.Range("A2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)

# 3. Turn the data range into an Excel Table ("Table1") spanning A1:U72
$tableRange = $ws.Range("A1:U72")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
